$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks that live on column B before deleting it, so no
# stray hyperlink definitions remain in the worksheet.
$ws.Hyperlinks.Delete()

# Delete column B ("Taxon name") entirely - this removes the hyperlinked
# plain-text species names column, shifting C:F left to B:E.
$ws.Columns.Item(2).Delete()

# The "Hyperlink" cell style (and its underline font) was only used by the
# now-deleted column B, so drop the now-unused named style too.
$wb.Styles.Item("Hyperlink").Delete()
